# F03 Froze Token Embeddings + Decoder 1
# Update the per-epoch accuracy values in column B (and refresh the
# "<__main__.DisplayOutputs object at ...>" memory-address labels in
# column A for the rows that carry them) to reflect the new training run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: epoch accuracy values -------------------------------------
$ws.Range("B3").Value = 0.96875
$ws.Range("B4:B5").Value = 0.953125
$ws.Range("B6:B7").Value = 0.890625
$ws.Range("B8").Value = 0.859375
$ws.Range("B9:B10").Value = 0.90625
$ws.Range("B11").Value = 0.890625
$ws.Range("B13").Value = 0.796875
$ws.Range("B14").Value = 0.828125
$ws.Range("B15").Value = 0.859375
$ws.Range("B16").Value = 0.875
$ws.Range("B17").Value = 0.8125
$ws.Range("B18").Value = 0.765625
$ws.Range("B19").Value = 0.859375
$ws.Range("B20").Value = 0.84375
$ws.Range("B21").Value = 0.8125
$ws.Range("B22:B23").Value = 0.734375
$ws.Range("B24").Value = 0.71875
$ws.Range("B25:B26").Value = 0.6875
$ws.Range("B27").Value = 0.671875
$ws.Range("B28:B29").Value = 0.6875
$ws.Range("B30:B40").Value = 0.671875
$ws.Range("B41:B64").Value = 0.6875
$ws.Range("B65:B102").Value = 0.703125
$ws.Range("B103").Value = 0.78125
$ws.Range("B104:B105").Value = 0.6875
$ws.Range("B106").Value = 0.703125
$ws.Range("B107").Value = 0.640625
$ws.Range("B108").Value = 0.671875
$ws.Range("B109").Value = 0.6875
$ws.Range("B110").Value = 0.671875
$ws.Range("B112").Value = 0.78125
$ws.Range("B113").Value = 0.578125
$ws.Range("B114").Value = 0.65625
$ws.Range("B115").Value = 0.796875
$ws.Range("B116").Value = 0.734375
$ws.Range("B117").Value = 0.78125
$ws.Range("B118").Value = 0.6885245901639344

# --- Column A: refresh the repr() memory address on the DisplayOutputs rows
$newRepr = "<__main__.DisplayOutputs object at 0x7fd8f0053640>"
$ws.Range("A102:A118").Value = $newRepr
